$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2593.6
$ws.Range("I38").Value = 2326.3333
$ws.Range("J38").Value = 4999
$ws.Range("K38").Value = 6978.999899999999
$ws.Range("L38").Value = 14997
$ws.Range("M38").Value = -6606.999899999999
$ws.Range("N38").Value = -15741

$ws.Range("H39").Value = 863.4286
$ws.Range("J39").Value = 1536.25
$ws.Range("L39").Value = 4608.75
$ws.Range("N39").Value = -5200.75

$ws.Range("H41").Value = 1141.3334
$ws.Range("J41").Value = 181.375
$ws.Range("L41").Value = 181.375
$ws.Range("N41").Value = -1061.375

$ws.Range("H42").Value = 73791.78999999999
$ws.Range("I42").Value = 79288.62
$ws.Range("J42").Value = 2333
$ws.Range("K42").Value = 237865.86
$ws.Range("L42").Value = 6999
$ws.Range("M42").Value = -237635.86
$ws.Range("N42").Value = -7459

$ws.Range("H80").Value = 701.0741
$ws.Range("I80").Value = 656.8125
$ws.Range("K80").Value = 1970.4375
$ws.Range("M80").Value = -972.4375

$ws.Range("H83").Value = 701.0741
$ws.Range("I83").Value = 656.8125
$ws.Range("K83").Value = 5911.3125
$ws.Range("M83").Value = -919.3125

$ws.Range("H99").Value = 2058.75
$ws.Range("I99").Value = 632.1818
$ws.Range("J99").Value = 3802.3333
$ws.Range("K99").Value = 1896.5454
$ws.Range("L99").Value = 11406.9999
$ws.Range("M99").Value = -398.5454
$ws.Range("N99").Value = -14402.9999

$ws.Range("H101").Value = 314.45456
$ws.Range("I101").Value = 312.42856
$ws.Range("J101").Value = 318
$ws.Range("K101").Value = 937.28568
$ws.Range("L101").Value = 954
$ws.Range("M101").Value = 684.71432
$ws.Range("N101").Value = -4198

$ws.Range("H132").Value = 19820.215
$ws.Range("I132").Value = 3731
$ws.Range("K132").Value = 11193
$ws.Range("M132").Value = -8663

$ws.Range("H138").Value = 3708.7273
$ws.Range("J138").Value = 4553.17
$ws.Range("L138").Value = 13659.51
$ws.Range("N138").Value = -23939.51


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 692.2
$ws.Range("I5").Value = 365.5
$ws.Range("J5").Value = 1999
$ws.Range("K5").Value = 365.5
$ws.Range("L5").Value = 1999
$ws.Range("M5").Value = -253.5
$ws.Range("N5").Value = -2223

$ws.Range("H10").Value = 6523.75
$ws.Range("J10").Value = 8665
$ws.Range("L10").Value = 8665
$ws.Range("N10").Value = -9005

$ws.Range("H32").Value = 2469.9385
$ws.Range("I32").Value = 1941.2
$ws.Range("J32").Value = 8814.799999999999
$ws.Range("K32").Value = 1941.2
$ws.Range("L32").Value = 8814.799999999999
$ws.Range("M32").Value = -1654.2
$ws.Range("N32").Value = -9388.799999999999

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 5669.717
$ws.Range("I132").Value = 2565.262
$ws.Range("K132").Value = 7695.786
$ws.Range("M132").Value = -5165.786


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 692.2
$ws.Range("I4").Value = 365.5
$ws.Range("J4").Value = 1999
$ws.Range("K4").Value = 365.5
$ws.Range("L4").Value = 1999
$ws.Range("M4").Value = -250.5
$ws.Range("N4").Value = -2229

$ws.Range("H17").Value = 600
$ws.Range("I17").Value = 600
$ws.Range("K17").Value = 600
$ws.Range("M17").Value = -428

$ws.Range("H20").Value = 3107.2856
$ws.Range("I20").Value = 2817
$ws.Range("K20").Value = 2817
$ws.Range("M20").Value = -2570

$ws.Range("H105").Value = 2050.625
$ws.Range("I105").Value = 1971.1364
$ws.Range("J105").Value = 2925
$ws.Range("K105").Value = 1971.1364
$ws.Range("L105").Value = 2925
$ws.Range("M105").Value = -224.1364000000001
$ws.Range("N105").Value = -6419

$ws.Range("H107").Value = 2971.0303
$ws.Range("I107").Value = 2113.92
$ws.Range("J107").Value = 5649.5
$ws.Range("K107").Value = 2113.92
$ws.Range("L107").Value = 5649.5
$ws.Range("M107").Value = -193.9200000000001
$ws.Range("N107").Value = -9489.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2500.9048
$ws.Range("I31").Value = 2097.8518
$ws.Range("K31").Value = 2097.8518
$ws.Range("M31").Value = -1802.8518

$ws.Range("H34").Value = 2500.9048
$ws.Range("I34").Value = 2097.8518
$ws.Range("K34").Value = 2097.8518
$ws.Range("M34").Value = -1895.8518

$ws.Range("H99").Value = 5718.7144
$ws.Range("I99").Value = 3396.3333
$ws.Range("J99").Value = 7460.5
$ws.Range("K99").Value = 3396.3333
$ws.Range("L99").Value = 7460.5
$ws.Range("M99").Value = -1898.3333
$ws.Range("N99").Value = -10456.5

$ws.Range("H122").Value = 3530.7576
$ws.Range("I122").Value = 3140.7778
$ws.Range("K122").Value = 9422.3334
$ws.Range("M122").Value = -6972.3334

$ws.Range("H126").Value = 5718.7144
$ws.Range("I126").Value = 3396.3333
$ws.Range("J126").Value = 7460.5
$ws.Range("K126").Value = 10188.9999
$ws.Range("L126").Value = 22381.5
$ws.Range("M126").Value = -7718.999899999999
$ws.Range("N126").Value = -27321.5

$ws.Range("H132").Value = 17929.484
$ws.Range("I132").Value = 19166
$ws.Range("J132").Value = 5564.3335
$ws.Range("K132").Value = 57498
$ws.Range("L132").Value = 16693.0005
$ws.Range("M132").Value = -54968
$ws.Range("N132").Value = -21753.0005


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 318.22223
$ws.Range("I16").Value = 363.33334
$ws.Range("K16").Value = 1090.00002
$ws.Range("M16").Value = -917.0000199999999

$ws.Range("H50").Value = 2306.75
$ws.Range("I50").Value = 132.14285
$ws.Range("J50").Value = 3998.111
$ws.Range("K50").Value = 396.42855
$ws.Range("L50").Value = 11994.333
$ws.Range("M50").Value = 84.57144999999997
$ws.Range("N50").Value = -12956.333

$ws.Range("H53").Value = 2306.75
$ws.Range("I53").Value = 132.14285
$ws.Range("J53").Value = 3998.111
$ws.Range("K53").Value = 396.42855
$ws.Range("L53").Value = 11994.333
$ws.Range("M53").Value = 84.57144999999997
$ws.Range("N53").Value = -12956.333

$ws.Range("H60").Value = 1450.5238
$ws.Range("I60").Value = 740.625
$ws.Range("J60").Value = 1887.3846
$ws.Range("K60").Value = 2221.875
$ws.Range("L60").Value = 5662.1538
$ws.Range("M60").Value = -1970.875
$ws.Range("N60").Value = -6164.1538

$ws.Range("H92").Value = 495.28125
$ws.Range("I92").Value = 386.125
$ws.Range("K92").Value = 1158.375
$ws.Range("M92").Value = 89.625


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26998.5
$ws.Range("I70").Value = 41142.145
$ws.Range("J70").Value = 12854.857
$ws.Range("K70").Value = 41142.145
$ws.Range("L70").Value = 12854.857
$ws.Range("M70").Value = -40872.145
$ws.Range("N70").Value = -13394.857

$ws.Range("H73").Value = 26998.5
$ws.Range("I73").Value = 41142.145
$ws.Range("J73").Value = 12854.857
$ws.Range("K73").Value = 41142.145
$ws.Range("L73").Value = 12854.857
$ws.Range("M73").Value = -40206.145
$ws.Range("N73").Value = -14726.857

$ws.Range("H122").Value = 2797.0625
$ws.Range("I122").Value = 2637.9167
$ws.Range("K122").Value = 7913.750100000001
$ws.Range("M122").Value = -5463.750100000001


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 43188.965
$ws.Range("I132").Value = 127430.5
$ws.Range("J132").Value = 9492.35
$ws.Range("K132").Value = 382291.5
$ws.Range("L132").Value = 28477.05
$ws.Range("M132").Value = -379761.5
$ws.Range("N132").Value = -33537.05


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 56750.188
$ws.Range("J46").Value = 56750.188
$ws.Range("L46").Value = 56750.188
$ws.Range("N46").Value = -57212.188

$ws.Range("H100").Value = 1093.1666
$ws.Range("I100").Value = 593.3
$ws.Range("J100").Value = 2092.9
$ws.Range("K100").Value = 1186.6
$ws.Range("L100").Value = 4185.8
$ws.Range("M100").Value = -645.5999999999999
$ws.Range("N100").Value = -5267.8

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H134").Value = 56750.188
$ws.Range("J134").Value = 56750.188
$ws.Range("L134").Value = 170250.564
$ws.Range("N134").Value = -175320.564

